$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($rowIndex, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowIndex, $i + 1).Value = $values[$i]
    }
}

# Insert a new weekly record (2022-03-23, $/paquete 36 unidades, Región Metropolitana)
# just above the existing 2022-03-11 entry -- pushes every row from the old row 4
# downward by one.
$ws.Rows.Item(4).Insert()
Set-RowValues 4 @(11, "Vega Monumental Concepción", "Bíobío", 44643, 8, 100112037, "Cebollín", "Sin especificar", "Primera", 180, 6500, 7000, 6778, "`$/paquete 36 unidades", "Región Metropolitana", 188, 36, "Hortaliza")

# Insert a second new weekly record (2022-03-17, $/paquete 36 unidades, Región
# Metropolitana) above what is now the 2022-02-03 entry (old row 25, now row 26
# after the first insert above).
$ws.Rows.Item(26).Insert()
Set-RowValues 26 @(11, "Vega Monumental Concepción", "Bíobío", 44637, 8, 100112037, "Cebollín", "Sin especificar", "Primera", 110, 6500, 7000, 6773, "`$/paquete 36 unidades", "Región Metropolitana", 188, 36, "Hortaliza")
